# Generate Report for Handoff
#
# The localization source file "0d44264d-ceb5-48ad-b8f0-a7dbd61ec138.md" was
# renamed/re-handed-off as "b4fe5b0b-4ce6-4a74-a8d0-dedd5801e610.md" (new
# handoff xlf + datetime), and a second source file
# "df1e1583-cfe0-47a3-9e6c-9d66b7db1622.md" showed up whose handoff
# transform failed. ".localization-config" keeps its row but moves down.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/b78ed5888ae1173f15d7719230486941ff017489"
$oldMd = "0d44264d-ceb5-48ad-b8f0-a7dbd61ec138.md"
$newMd = "b4fe5b0b-4ce6-4a74-a8d0-dedd5801e610.md"
$failedMd = "df1e1583-cfe0-47a3-9e6c-9d66b7db1622.md"
$cfgFile = ".localization-config"

$oldXlfBaseZh = "0d44264d-ceb5-48ad-b8f0-a7dbd61ec138.a04bc578a144bf95fdac0f1649eb825ec7417ac7.zh-cn.xlf"
$newXlfBaseZh = "b4fe5b0b-4ce6-4a74-a8d0-dedd5801e610.529a8dda16127c2de88b6d910597100f0baa5b10.zh-cn.xlf"
$oldXlfBaseDe = "0d44264d-ceb5-48ad-b8f0-a7dbd61ec138.a04bc578a144bf95fdac0f1649eb825ec7417ac7.de-de.xlf"
$newXlfBaseDe = "b4fe5b0b-4ce6-4a74-a8d0-dedd5801e610.529a8dda16127c2de88b6d910597100f0baa5b10.de-de.xlf"

$handoffZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c6d587191306d147e0f87bae80584c33bd25fa69/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/tianzh/$newXlfBaseZh"
$handoffDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40714a83b28c6f01d888cc383a58c008b6069b22/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/tianzh/$newXlfBaseDe"

$newHandoffDtZh = "2016-01-27 08:16:16"
$newHandoffDtDe = "2016-01-27 08:16:28"
$neverDt = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop every hyperlink on the sheet so we can rebuild them against the
# correct cells/targets without leaving stale duplicate entries behind.
$wsOverview.Range("A1").Hyperlinks.Delete()

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "$repoBase/e2e/$newMd", "", "", $newMd)

$wsOverview.Range("B3").Value = "Handoff transform failed"
$wsOverview.Range("C3").Value = "Handoff transform failed"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "$repoBase/e2e/$failedMd", "", "", $failedMd)

$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$repoBase/$cfgFile", "", "", $cfgFile)

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A1").Hyperlinks.Delete()

# Row 2 -- renamed source file, refreshed handoff artifact + timestamp.
$wsZh.Range("B2").Value = "Ready for handoff"
$wsZh.Range("D2").Value = $newHandoffDtZh
$wsZh.Range("G2").Value = $neverDt
$wsZh.Range("H2").Value = "Include"
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoBase/e2e/$newMd", "", "", $newMd)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $handoffZhUrl, "", "", $newXlfBaseZh)

# Row 3 (new) -- second source file, handoff transform failed.
$wsZh.Range("B3").Value = "Handoff transform failed"
$wsZh.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("D3").Value = $neverDt
$wsZh.Range("G3").Value = $neverDt
$wsZh.Range("H3").Value = "Ignored"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$repoBase/e2e/$failedMd", "", "", $failedMd)

# Row 4 (new) -- .localization-config, shifted down from row 3.
$wsZh.Range("B4").Value = "Not to be localized"
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("D4").Value = $neverDt
$wsZh.Range("G4").Value = $neverDt
$wsZh.Range("H4").Value = "Ignored"
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$repoBase/$cfgFile", "", "", $cfgFile)

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as zh-cn
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A1").Hyperlinks.Delete()

# Row 2 -- renamed source file, refreshed handoff artifact + timestamp.
$wsDe.Range("B2").Value = "Ready for handoff"
$wsDe.Range("D2").Value = $newHandoffDtDe
$wsDe.Range("G2").Value = $neverDt
$wsDe.Range("H2").Value = "Include"
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBase/e2e/$newMd", "", "", $newMd)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $handoffDeUrl, "", "", $newXlfBaseDe)

# Row 3 (new) -- second source file, handoff transform failed.
$wsDe.Range("B3").Value = "Handoff transform failed"
$wsDe.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("D3").Value = $neverDt
$wsDe.Range("G3").Value = $neverDt
$wsDe.Range("H3").Value = "Ignored"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoBase/e2e/$failedMd", "", "", $failedMd)

# Row 4 (new) -- .localization-config, shifted down from row 3.
$wsDe.Range("B4").Value = "Not to be localized"
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("D4").Value = $neverDt
$wsDe.Range("G4").Value = $neverDt
$wsDe.Range("H4").Value = "Ignored"
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$repoBase/$cfgFile", "", "", $cfgFile)
